# Natmi following Dr Hou advice
# Updates the recomputed NATMI LR-pair statistics for Vegfa-Gpc1 (rows 2-17)
# following re-analysis with 3 ligand-/receptor-expressing cells (columns E, K)
# instead of 1, and the consequently recalculated derived metrics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> hashtable of column letter -> new value
$updates = @{
    2  = @{ E=3; G=28.479168; H=85.437504; I=0.4446244458164738; J=0.4446244458164738; K=3; M=0.8478306666666667; N=2.543492; O=0.01460351867535248; P=0.01460351867535248; Q=24.145511991552; R=217.309607923968; S=0.006493081397999124; T=0.006493081397999123 }
    3  = @{ E=3; G=28.479168; H=85.437504; I=0.4446244458164738; J=0.4446244458164738; K=3; M=7.020353; N=21.061059; O=0.1209225617494376; P=0.1209225617494376; Q=199.933812506304; R=1799.404312556736; S=0.05376512700455201; T=0.05376512700455201 }
    4  = @{ E=3; G=28.479168; H=85.437504; I=0.4446244458164738; J=0.4446244458164738; K=3; M=1.462291666666667; N=4.386875; O=0.02518734518879435; P=0.02518734518879435; Q=41.64485004; R=374.80365036; S=0.01119890939615592; T=0.01119890939615592 }
    5  = @{ E=3; G=28.479168; H=85.437504; I=0.4446244458164738; J=0.4446244458164738; K=3; M=48.72612633333333; N=146.178379; O=0.8392865743864156; P=0.8392865743864156; Q=1387.679537836224; R=12489.11584052602; S=0.3731673280177668; T=0.3731673280177668 }
    6  = @{ E=3; G=18.12667766666667; H=54.380033; I=0.2829985767855128; J=0.2829985767855128; K=3; M=0.8478306666666667; N=2.543492; O=0.01460351867535248; P=0.01460351867535248; Q=15.36835321058178; R=138.315178895236; S=0.00413277500118541; T=0.004132775001185409 }
    7  = @{ E=3; G=18.12667766666667; H=54.380033; I=0.2829985767855128; J=0.2829985767855128; K=3; M=7.020353; N=21.061059; O=0.1209225617494376; P=0.1209225617494376; Q=127.2556759372164; R=1145.301083434947; S=0.03422091287634912; T=0.03422091287634912 }
    8  = @{ E=3; G=18.12667766666667; H=54.380033; I=0.2829985767855128; J=0.2829985767855128; K=3; M=1.462291666666667; N=4.386875; O=0.02518734518879435; P=0.02518734518879435; Q=26.50648969631944; R=238.558407266875; S=0.007127982841434235; T=0.007127982841434235 }
    9  = @{ E=3; G=18.12667766666667; H=54.380033; I=0.2829985767855128; J=0.2829985767855128; K=3; M=48.72612633333333; N=146.178379; O=0.8392865743864156; P=0.8392865743864156; Q=883.242785989612; R=7949.185073906508; S=0.2375169060665441; T=0.2375169060665441 }
    10 = @{ E=3; G=11.513346; H=34.540038; I=0.179749460544048; J=0.179749460544048; K=3; M=0.8478306666666667; N=2.543492; O=0.01460351867535248; P=0.01460351867535248; Q=9.761367814744; R=87.85231033269601; S=0.00262497460393954; T=0.002624974603939539 }
    11 = @{ E=3; G=11.513346; H=34.540038; I=0.179749460544048; J=0.179749460544048; K=3; M=7.020353; N=21.061059; O=0.1209225617494376; P=0.1209225617494376; Q=80.827753131138; R=727.4497781802421; S=0.02173576524206574; T=0.02173576524206574 }
    12 = @{ E=3; G=11.513346; H=34.540038; I=0.179749460544048; J=0.179749460544048; K=3; M=1.462291666666667; N=4.386875; O=0.02518734518879435; P=0.02518734518879435; Q=16.83586991125; R=151.52282920125; S=0.004527411710222508; T=0.004527411710222508 }
    13 = @{ E=3; G=11.513346; H=34.540038; I=0.179749460544048; J=0.179749460544048; K=3; M=48.72612633333333; N=146.178379; O=0.8392865743864156; P=0.8392865743864156; Q=561.000751715378; R=5049.006765438403; S=0.1508613089878202; T=0.1508613089878202 }
    14 = @{ E=3; G=5.932994999999999; H=17.798985; I=0.09262751685396531; J=0.09262751685396531; K=3; M=0.8478306666666667; N=2.543492; O=0.01460351867535248; P=0.01460351867535248; Q=5.03017510618; R=45.27157595561999; S=0.001352687672228409; T=0.001352687672228409 }
    15 = @{ E=3; G=5.932994999999999; H=17.798985; I=0.09262751685396531; J=0.09262751685396531; K=3; M=7.020353; N=21.061059; O=0.1209225617494376; P=0.1209225617494376; Q=41.65171924723499; R=374.865473225115; S=0.01120075662647069; T=0.01120075662647069 }
    16 = @{ E=3; G=5.932994999999999; H=17.798985; I=0.09262751685396531; J=0.09262751685396531; K=3; M=1.462291666666667; N=4.386875; O=0.02518734518879435; P=0.02518734518879435; Q=8.675769146874998; R=78.08192232187498; S=0.00233304124098169; T=0.00233304124098169 }
    17 = @{ E=3; G=5.932994999999999; H=17.798985; I=0.09262751685396531; J=0.09262751685396531; K=3; M=48.72612633333333; N=146.178379; O=0.8392865743864156; P=0.8392865743864156; Q=289.091863905035; R=2601.826775145315; S=0.07774103131428452; T=0.07774103131428452 }
}

foreach ($rowNum in $updates.Keys) {
    $cols = $updates[$rowNum]
    foreach ($colLetter in $cols.Keys) {
        $ws.Range("$colLetter$rowNum").Value = $cols[$colLetter]
    }
}
